# Adapt illustrative case to reduce min up-/down-time
# Updates the "ScenarioA" sheet, row 8 (the single thermal generator entry):
#   MinProd (H8)    200 -> 100
#   MinUpTime (K8)    5 -> 1
#   MinDownTime (L8)  5 -> 1
# Also restores the window/selection state recorded by Excel after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioA")

$ws.Range("H8").Value = 100
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 1

# Move the active selection to K8, matching the post-edit cursor position.
$ws.Range("K8").Select()

# Restore the workbook window geometry recorded on save.
$excel.ActiveWindow.Top = -11250
$excel.ActiveWindow.Height = 10545
$excel.ActiveWindow.Width = 19410
$excel.ActiveWindow.Left = 50205
